$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.51"
$ws.Range("E2").Value = "'1.39%"
$ws.Range("D3").Value = "'32.77"
$ws.Range("E3").Value = "'4.69%"
$ws.Range("D4").Value = "'4.945"
$ws.Range("E4").Value = "'-3.13%"
$ws.Range("D5").Value = "'0.07847"
$ws.Range("E5").Value = "'-1.32%"
$ws.Range("D6").Value = "'2.029"
$ws.Range("E6").Value = "'-11.07%"
$ws.Range("D7").Value = "'7.837"
$ws.Range("E7").Value = "'0.85%"
$ws.Range("D8").Value = "'3.809"
$ws.Range("E8").Value = "'-1.46%"
$ws.Range("D9").Value = "'0.9239"
$ws.Range("E9").Value = "'-0.24%"
$ws.Range("D10").Value = "'0.1755"
$ws.Range("E10").Value = "'1.11%"
$ws.Range("D11").Value = "'0.07832"
$ws.Range("E11").Value = "'3.87%"
$ws.Range("D12").Value = "'0.08673"
$ws.Range("E12").Value = "'-7.33%"
$ws.Range("D13").Value = "'0.03143"
$ws.Range("E13").Value = "'3.18%"
$ws.Range("D14").Value = "'0.1006"
$ws.Range("D15").Value = "'0.001519"
$ws.Range("E15").Value = "'0.59%"
$ws.Range("D16").Value = "'0.005860"
$ws.Range("E16").Value = "'-1.00%"
$ws.Range("E17").Value = "'-0.53%"
$ws.Range("D18").Value = "'2.155"
$ws.Range("E18").Value = "'-4.98%"
$ws.Range("D19").Value = "'0.3307"
$ws.Range("E20").Value = "'-3.38%"
$ws.Range("D21").Value = "'4.319"
$ws.Range("E21").Value = "'10.06%"
$ws.Range("E22").Value = "'16.99%"
$ws.Range("D23").Value = "'0.04564"
$ws.Range("E23").Value = "'-0.99%"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").Value = "'-1.98%"
$ws.Range("D25").Value = "'0.004444"
$ws.Range("E25").Value = "'-0.85%"
$ws.Range("E26").Value = "'4.31%"
$ws.Range("D39").Value = "'0.01738"
$ws.Range("E39").Value = "'-1.04%"
$ws.Range("D40").Value = "'0.04793"
$ws.Range("E40").Value = "'3.56%"
$ws.Range("D41").Value = "'0.007475"
$ws.Range("E41").Value = "'7.43%"
$ws.Range("D42").Value = "'0.1361"
$ws.Range("E42").Value = "'-0.11%"
$ws.Range("D43").Value = "'0.002341"
$ws.Range("E43").Value = "'7.01%"
$ws.Range("D44").Value = "'0.01058"
$ws.Range("E44").Value = "'3.27%"
$ws.Range("D45").Value = "'0.00006258"
$ws.Range("E45").Value = "'-0.26%"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("E47").Value = "'-61.14%"
$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'-29.08%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.04%"
